$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "24.03.2023 16:20 (CET)"
$ws.Range("C5").Value = "https://gitlab.intra.infineon.com/semantic-web-projects/digital-reference/process_model_version/-/commit/7d95b41a9b432c782ac3bd45ad803e677399ab28"
$ws.Range("D5").Value = "b4ab941a87685341ec282f6e0ff6634ae8a1afa78e3fffdcbad16e151dde0c33"
